$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells (prices in column D look numeric, e.g. "1.00",
# "0.0790", "311.17") stay stored as text instead of being auto-coerced
# into numbers by Excel, matching the original inline-string formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.806.02"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.527.59"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.17"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.25"
$ws.Range("E6").Value = "  +2.66%  "
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -1.74%  "
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.917.55"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.570.57"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.32"
$ws.Range("E16").Value = "  -3.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.815"
$ws.Range("E17").Value = "  -3.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.774.10"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.83"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.68"
$ws.Range("E23").Value = "  -1.84%  "
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.49"
$ws.Range("E27").Value = "  -5.61%  "
$ws.Range("E28").Value = "  -2.71%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("E30").Value = "  -3.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "161.81"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.82"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.80"
$ws.Range("E33").Value = "  +8.81%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0790"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.30"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.08"
$ws.Range("E37").Value = "  -6.08%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.96"
$ws.Range("E38").Value = "  -7.18%  "
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.17"
$ws.Range("E41").Value = "  +1.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.08"
$ws.Range("E42").Value = "  -3.13%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("E44").Value = "  +3.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0299"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.991.05"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.05"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.771.00"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("E49").Value = "  -1.77%  "
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.48"
$ws.Range("E51").Value = "  -1.97%  "
